$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows to append: row 130 ("四方坪站充电量(kw)") and row 131 ("高岭站充电量(kw)")
# Column A: date serial 45965 (2025-11-04)
# Column B: station label (same text as rows 128/129 -> reuse their shared string)
# Columns C:Z: hourly charge values

$row130 = @{
    A = 45965
    C = 580.29999999999995
    D = 1302.6900000000003
    E = 512.1869999999999
    F = 497.57899999999995
    G = 188.35500000000002
    H = 647.16300000000012
    I = 621.56000000000017
    J = 106.857
    K = 121.711
    L = 127.67999999999999
    M = 259.02499999999998
    N = 356.32799999999997
    O = 902.93599999999992
    P = 1428.2729999999997
    Q = 525.58699999999999
    R = 362.791
    S = 232.25
    T = 135.70700000000002
    U = 60.6
    V = 76.960000000000008
    W = 50.555
    X = 160.46
    Y = 0
    Z = 28.19
}

$row131 = @{
    A = 45965
    C = 346.10400000000004
    D = 197.57
    E = 124.304
    F = 79.046999999999997
    G = 43.555
    H = 178.89699999999999
    I = 332.29399999999998
    J = 116.553
    K = 377.46199999999999
    L = 153.18100000000001
    M = 70.730999999999995
    N = 272.697
    O = 411.57400000000001
    P = 455.19599999999997
    Q = 250.80699999999996
    R = 183.779
    S = 85.147000000000006
    T = 54.954000000000008
    U = 21.207000000000001
    V = 0
    W = 81.977999999999994
    X = 0
    Y = 0
    Z = 0
}

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")

# Copy rows 128:129 (values + formats, including the shared-string text in B
# and the date/number styles) straight down into the new rows 130:131 -
# mirrors dragging the fill handle / copy-paste in the UI.
$ws.Range("A128:Z129").Copy() | Out-Null
$ws.Range("A130").PasteSpecial(-4104) | Out-Null  # xlPasteAll
$excel.CutCopyMode = 0

# Now overwrite the numeric cells (everything except column B, whose text
# stays identical to rows 128/129) with the real data for the new date.
foreach ($col in $cols) {
    if ($col -eq "B") { continue }
    $addr130 = $col + "130"
    $ws.Range($addr130).Value = $row130[$col]
}
foreach ($col in $cols) {
    if ($col -eq "B") { continue }
    $addr131 = $col + "131"
    $ws.Range($addr131).Value = $row131[$col]
}
